$wb = $excel.ActiveWorkbook

# 1. Insert a new worksheet "Tabelle1" right after "Algorithm" (before "Operators")
$algSheet = $wb.Worksheets.Item("Algorithm")
$newSheet = $wb.Worksheets.Add($null, $algSheet)
$newSheet.Name = "Tabelle1"
$newSheet.Range("A1").Value = "select"
$newSheet.Range("B1").Value = "tools"
$newSheet.Range("C1").Value = "selTournamentDCD"
$newSheet.PageSetup.TopMargin = 56.692913399999995
$newSheet.PageSetup.BottomMargin = 56.692913399999995

# 2. On "Algorithm" sheet: update C4 value, selection, column width
$algSheet.Range("C4").Value = "selNSGA2revised"
$algSheet.Columns("B").ColumnWidth = 22.7109375
$algSheet.Range("C4").Select()

# 3. On "Operators" sheet: update C1 value, selection; make it the active/selected tab
$opSheet = $wb.Worksheets.Item("Operators")
$opSheet.Range("C1").Value = "UFTournSelection"
$opSheet.Range("C1").Select()
$opSheet.Select()
